# Adds a "FineshedAt" timestamp column to each of the data-writing
# worksheets (mirrors the Java test fixture being refactored to log a
# finish time per web-table write). The header cell gets the same
# highlighted fill style already used for the other header cells on
# that sheet, and only row 4 (the 3rd data row) receives a logged
# timestamp value, matching the original commit.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function Add-FineshedAtColumn {
    param(
        [string]$SheetName,
        [string]$HeaderCellRef,
        [string]$HeaderStyleSourceRef,
        [string]$TimestampCellRef,
        [string]$TimestampValue
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # New header cell text.
    $ws.Range($HeaderCellRef).Value = "FineshedAt"

    # Copy the look of an existing header cell (the highlighted fill)
    # onto the new header cell.
    $ws.Range($HeaderStyleSourceRef).Copy() | Out-Null
    $ws.Range($HeaderCellRef).PasteSpecial($xlPasteFormats) | Out-Null

    # Logged finish timestamp for this sheet's write operation.
    $ws.Range($TimestampCellRef).Value = $TimestampValue
}

Add-FineshedAtColumn "writeCompanyColumnIntoXcel" "B1" "A1" "B4" `
    "1571334563820`nThu Oct 17 10:49:23 PDT 2019"

Add-FineshedAtColumn "writeContactColumnIntoXcel" "B1" "A1" "B4" `
    "1571334565652`nThu Oct 17 10:49:25 PDT 2019"

Add-FineshedAtColumn "writeCountryColumnIntoXcel" "B1" "A1" "B4" `
    "1571334566494`nThu Oct 17 10:49:26 PDT 2019"

Add-FineshedAtColumn "writeWholeTableNestedFor" "D1" "A1" "D4" `
    "1571334569471`nThu Oct 17 10:49:29 PDT 2019"

Add-FineshedAtColumn "writeWholeTableSingleFor" "D1" "A1" "D4" `
    "1571334570770`nThu Oct 17 10:49:30 PDT 2019"
